$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 44, pushing the existing rows 44-59 down to 46-61.
$ws.Rows.Item(44).Insert()
$ws.Rows.Item(44).Insert()

# New row 44: Caigua, "Primera", week of 2021-09-29 (serial 44468)
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").Value = 44468
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 100112036
$ws.Range("G44").Value = "Caigua"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 120
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 8000
$ws.Range("M44").Value = 7500
$ws.Range("N44").Value = "`$/caja 20 kilos"
$ws.Range("O44").Value = "Región de Arica y Parinacota"
$ws.Range("P44").Value = 375
$ws.Range("Q44").Value = 20
$ws.Range("R44").Value = "Hortaliza"

# New row 45: Caigua, "Segunda", week of 2021-09-29 (serial 44468)
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 44468
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112036
$ws.Range("G45").Value = "Caigua"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Segunda"
$ws.Range("J45").Value = 120
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = 6500
$ws.Range("N45").Value = "`$/caja 20 kilos"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 325
$ws.Range("Q45").Value = 20
$ws.Range("R45").Value = "Hortaliza"
